$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 793
$ws.Range("I28").Value = 145
$ws.Range("J28").Value = 2575
$ws.Range("K28").Value = 145
$ws.Range("L28").Value = 2575
$ws.Range("M28").Value = 340
$ws.Range("N28").Value = -3545

$ws.Range("H31").Value = 774.5
$ws.Range("I31").Value = 774.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2323.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2093.5

$ws.Range("H34").Value = 1499.75
$ws.Range("I34").Value = 1499.75
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1499.75
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1296.75
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 1499.75
$ws.Range("I36").Value = 1499.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1499.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -784.75
$ws.Range("N36").ClearContents()

$ws.Range("H40").Value = 3610
$ws.Range("I40").Value = 1140
$ws.Range("J40").Value = 4433.3335
$ws.Range("K40").Value = 1140
$ws.Range("L40").Value = 4433.3335
$ws.Range("M40").Value = -965
$ws.Range("N40").Value = -4783.3335

$ws.Range("H43").Value = 1783.421
$ws.Range("I43").Value = 1100.3334
$ws.Range("J43").Value = 1911.5
$ws.Range("K43").Value = 1100.3334
$ws.Range("L43").Value = 1911.5
$ws.Range("M43").Value = -1031.3334
$ws.Range("N43").Value = -2049.5

$ws.Range("H70").Value = 2828.0386
$ws.Range("I70").Value = 1326.6666
$ws.Range("J70").Value = 3023.8696
$ws.Range("K70").Value = 3979.9998
$ws.Range("L70").Value = 9071.6088
$ws.Range("M70").Value = -3709.9998
$ws.Range("N70").Value = -9611.6088

$ws.Range("H73").Value = 2828.0386
$ws.Range("I73").Value = 1326.6666
$ws.Range("J73").Value = 3023.8696
$ws.Range("K73").Value = 3979.9998
$ws.Range("L73").Value = 9071.6088
$ws.Range("M73").Value = -3043.9998
$ws.Range("N73").Value = -10943.6088

$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 5000
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -6872

$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 25000
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -34360

$ws.Range("H112").Value = 4311520
$ws.Range("I112").Value = 2048.5715
$ws.Range("J112").Value = 4903016
$ws.Range("K112").Value = 6145.7145
$ws.Range("L112").Value = 14709048
$ws.Range("M112").Value = -5037.7145
$ws.Range("N112").Value = -14711264

$ws.Range("H132").Value = 4001824.5
$ws.Range("I132").Value = 4652590
$ws.Range("J132").Value = 4264.2856
$ws.Range("K132").Value = 13957770
$ws.Range("L132").Value = 12792.8568
$ws.Range("M132").Value = -13955240
$ws.Range("N132").Value = -17852.8568

$ws.Range("H137").Value = 4426.095
$ws.Range("I137").Value = 4396.75
$ws.Range("J137").Value = 4520
$ws.Range("K137").Value = 13190.25
$ws.Range("L137").Value = 13560
$ws.Range("M137").Value = -10640.25
$ws.Range("N137").Value = -18660

$ws.Range("H138").Value = 3680.3699
$ws.Range("I138").Value = 1626.5
$ws.Range("J138").Value = 7179.5557
$ws.Range("K138").Value = 4879.5
$ws.Range("L138").Value = 21538.6671
$ws.Range("M138").Value = 260.5
$ws.Range("N138").Value = -31818.6671

$ws.Range("H141").Value = 589737.8
$ws.Range("I141").Value = 1625.9375
$ws.Range("J141").Value = 2471695.8
$ws.Range("K141").Value = 4877.8125
$ws.Range("L141").Value = 7415087.399999999
$ws.Range("M141").Value = 302.1875
$ws.Range("N141").Value = -7425447.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4072
$ws.Range("I32").Value = 2789.9
$ws.Range("J32").Value = 14043.889
$ws.Range("K32").Value = 2789.9
$ws.Range("L32").Value = 14043.889
$ws.Range("M32").Value = -2502.9
$ws.Range("N32").Value = -14617.889

$ws.Range("H61").Value = 2454.86
$ws.Range("I61").Value = 923.6667
$ws.Range("J61").Value = 3563.6553
$ws.Range("K61").Value = 923.6667
$ws.Range("L61").Value = 3563.6553
$ws.Range("M61").Value = -711.6667

$ws.Range("H74").Value = 742.2857
$ws.Range("I74").Value = 418.66666
$ws.Range("J74").Value = 1173.7778
$ws.Range("K74").Value = 418.66666
$ws.Range("L74").Value = 1173.7778
$ws.Range("M74").Value = 455.33334
$ws.Range("N74").Value = -2921.7778

$ws.Range("H77").Value = 742.2857
$ws.Range("I77").Value = 418.66666
$ws.Range("J77").Value = 1173.7778
$ws.Range("K77").Value = 2093.3333
$ws.Range("L77").Value = 5868.889
$ws.Range("M77").Value = 2274.6667
$ws.Range("N77").Value = -14604.889

$ws.Range("H132").Value = 19611208
$ws.Range("I132").Value = 28574910
$ws.Range("J132").Value = 3106.25
$ws.Range("K132").Value = 85724730
$ws.Range("L132").Value = 9318.75
$ws.Range("M132").Value = -85722200
$ws.Range("N132").Value = -14378.75

$ws.Range("H136").Value = 2454.86
$ws.Range("I136").Value = 923.6667
$ws.Range("J136").Value = 3563.6553
$ws.Range("K136").Value = 2771.0001
$ws.Range("L136").Value = 10690.9659
$ws.Range("M136").Value = -221.0001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3666.6667
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 500
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 1420
$ws.Range("N107").Value = -13840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2727.6667
$ws.Range("I99").Value = 1400
$ws.Range("J99").Value = 2893.625
$ws.Range("K99").Value = 1400
$ws.Range("L99").Value = 2893.625
$ws.Range("M99").Value = 98
$ws.Range("N99").Value = -5889.625

$ws.Range("H126").Value = 2727.6667
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 2893.625
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 8680.875
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -13620.875

$ws.Range("H134").Value = 1864.3889
$ws.Range("I134").Value = 1445.122
$ws.Range("J134").Value = 3186.6924
$ws.Range("K134").Value = 4335.366
$ws.Range("L134").Value = 9560.0772
$ws.Range("M134").Value = -1800.366

$ws.Range("H135").Value = 26181.818
$ws.Range("I135").Value = 30000
$ws.Range("J135").Value = 25800
$ws.Range("K135").Value = 30000
$ws.Range("L135").Value = 25800
$ws.Range("M135").Value = -24930
$ws.Range("N135").Value = -35940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6254735
$ws.Range("I139").Value = 7577613.5
$ws.Range("J139").Value = 18307.572
$ws.Range("K139").Value = 22732840.5
$ws.Range("L139").Value = 54922.716
$ws.Range("M139").Value = -22727700.5
$ws.Range("N139").Value = -65202.716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2547.5
$ws.Range("I43").Value = 2547.5
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 2547.5
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -2396.5
$ws.Range("N43").ClearContents()

$ws.Range("H132").Value = 5156
$ws.Range("I132").Value = 5004.8
$ws.Range("J132").Value = 5250.5
$ws.Range("K132").Value = 15014.4
$ws.Range("L132").Value = 15751.5
$ws.Range("M132").Value = -12484.4
$ws.Range("N132").Value = -20811.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 3017
$ws.Range("I29").Value = 3017
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 3017
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -2722

$ws.Range("H68").Value = 1708
$ws.Range("I68").Value = 1031.8182
$ws.Range("J68").Value = 6666.6665
$ws.Range("K68").Value = 1031.8182
$ws.Range("L68").Value = 6666.6665
$ws.Range("M68").Value = -282.8181999999999

$ws.Range("H71").Value = 1708
$ws.Range("I71").Value = 1031.8182
$ws.Range("J71").Value = 6666.6665
$ws.Range("K71").Value = 5159.090999999999
$ws.Range("L71").Value = 33333.3325
$ws.Range("M71").Value = -1415.090999999999

$ws.Range("H93").Value = 1144.0869
$ws.Range("I93").Value = 794.625
$ws.Range("J93").Value = 1942.8572
$ws.Range("K93").Value = 794.625
$ws.Range("L93").Value = 1942.8572
$ws.Range("M93").Value = 453.375
$ws.Range("N93").Value = -4438.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 70019
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 70019
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 70019
$ws.Range("N28").Value = -70715

$ws.Range("H39").Value = 17864.143
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 17864.143
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 17864.143
$ws.Range("N39").Value = -18690.143

$ws.Range("H43").Value = 3833.3333
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -1351
$ws.Range("N43").Value = -5298

$ws.Range("H46").Value = 47000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 47000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 47000
$ws.Range("N46").Value = -47462

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H82").Value = 28529.629
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 28529.629
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 28529.629
$ws.Range("N82").Value = -29295.629

$ws.Range("H85").Value = 28529.629
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 28529.629
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 28529.629
$ws.Range("N85").Value = -31181.629

$ws.Range("H132").Value = 7488.9614
$ws.Range("I132").Value = 3147.2
$ws.Range("J132").Value = 13409.546
$ws.Range("K132").Value = 9441.599999999999
$ws.Range("L132").Value = 40228.638
$ws.Range("M132").Value = -6911.599999999999
$ws.Range("N132").Value = -45288.638

$ws.Range("H134").Value = 47000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 47000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 141000
$ws.Range("N134").Value = -146070
